$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (between "Partner Code (*)" and "Charge Code (*)")
# This shifts all columns from E onward one position to the right.
$ws.Columns("E:E").Insert()

# Populate the new column's header and the sample data row with the new
# "OBH Partner" field.
$ws.Range("E1").Value2 = "OBH Partner"
$ws.Range("E2").Value2 = "26784508"

# Update the hidden _FilterDatabase defined name so it spans the new last
# column (R instead of Q).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$R`$1"
    }
}
